$wb = $excel.ActiveWorkbook

# --- Arkusz3 (sheet3.xml) updates -----------------------------------------
$ws3 = $wb.Worksheets.Item("Arkusz3")

# Header: "stacje" -> "Liczba"
$ws3.Range("B1").Value = "Liczba"

# Gdansk station count update
$ws3.Range("B3").Value = 121

# Re-labelled / re-ordered cities with their station counts
$ws3.Range("A6").Value = "Poznań"
$ws3.Range("B6").Value = 65

$ws3.Range("A7").Value = "Wrocław"
$ws3.Range("B7").Value = 55

$ws3.Range("B8").Value = 54
$ws3.Range("B9").Value = 43
$ws3.Range("B10").Value = 42

$ws3.Range("A11").Value = "Szczecin"
$ws3.Range("B11").Value = 42

$ws3.Range("B12").Value = 38
$ws3.Range("B13").Value = 37
$ws3.Range("B14").Value = 36
$ws3.Range("B15").Value = 33
$ws3.Range("B16").Value = 33

# Make Arkusz3 the active sheet / selection, which also clears the
# previously active Arkusz6 tab-selected flag.
$ws3.Activate()
$ws3.Range("B1").Select()
